$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

# --- 1. "Textfeld 22" (id 23): Abstandssensor IN -> Abstandssensor Echo IN ---
$txtEchoIn = $s.Shapes.Item(27)
$txtEchoIn.TextFrame.TextRange.Text = "Abstandssensor Echo IN"

# --- 2. "Pfeil: nach rechts 24" (id 25): rotate 180 degrees (rot="10800000") ---
$arrowTriggerOut = $s.Shapes.Item(28)
$arrowTriggerOut.Rotation = 180

# --- 3. "Textfeld 32" (id 33): reposition/resize + Abstandssensor IN -> Abstandssensor Trigger OUT ---
$txtTriggerOut = $s.Shapes.Item(29)
$txtTriggerOut.Left = 6.500039577484131
$txtTriggerOut.Width = 153.9864959716797
$txtTriggerOut.TextFrame.TextRange.Text = "Abstandssensor Trigger OUT"

# --- 4. New shape: "Pfeil: nach rechts 5" (id 6), duplicated from the just-rotated arrow ---
$newArrowRange = $arrowTriggerOut.Duplicate()
$newArrow = $newArrowRange.Item(1)
$newArrow.Name = "Pfeil: nach rechts 5"
$newArrow.Left = 160.48696899414062
$newArrow.Top = 245.1551513671875

# --- 5. New shape: "Textfeld 6" (id 7), duplicated from the just-edited text label ---
$newTextRange = $txtTriggerOut.Duplicate()
$newText = $newTextRange.Item(1)
$newText.Name = "Textfeld 6"
$newText.Left = 6.500511646270752
$newText.Top = 237.81515502929688
$newText.TextFrame.TextRange.Text = "Lautsprecher"
